# Invesco India Largecap Fund - insert an "Industry" column (new column C)
# between "Stock Name" and "Mutual Fund", shifting the existing
# Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ columns one to the
# right, and populate the new column with each holding's industry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:I -> D:J (preserves all existing values/styles/types).
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 3).Value = "Industry"

$industries = @(
    "Banks",
    "Banks",
    "Capital Markets",
    "IT - Software",
    "Construction",
    "Telecom - Services",
    "Finance",
    "IT - Software",
    "Transport Services",
    "Retailing",
    "Pharmaceuticals & Biotechnology",
    "Automobiles",
    "Consumer Durables",
    "Consumer Durables",
    "Cement & Cement Products",
    "Finance",
    "Healthcare Services",
    "Capital Markets",
    "Banks",
    "Electrical Equipment",
    "Automobiles",
    "Automobiles",
    "Capital Markets",
    "Healthcare Equipment & Supplies",
    "Banks",
    "Electrical Equipment",
    "Aerospace & Defense",
    "IT - Software",
    "Industrial Products",
    "Retailing",
    "Pharmaceuticals & Biotechnology",
    "Finance",
    "Diversified FMCG",
    "Non - Ferrous Metals",
    "Auto Components",
    "Healthcare Services",
    "Chemicals & Petrochemicals",
    "Financial Technology (Fintech)",
    "Finance",
    "Auto Components",
    "Healthcare Services",
    "Transport Services",
    "Realty",
    "Aerospace & Defense",
    "Consumer Durables",
    "Power",
    "Realty",
    "Retailing",
    "Industrial Manufacturing",
    "Food Products",
    "Banks",
    "Auto Components",
    "Capital Markets",
    "Petroleum Products",
    "Ferrous Metals",
    "Automobiles",
    "Pharmaceuticals & Biotechnology",
    "Food Products",
    "Banks"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
